$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-11-07 Thursday" "2024-11-08 Friday"

Replace-Text "839÷7=" "266÷3="
Replace-Text "992÷6=" "220÷5="
Replace-Text "109÷5=" "483÷4="
Replace-Text "487÷7=" "940÷9="
Replace-Text "814÷6=" "923÷4="
Replace-Text "214÷9=" "560÷7="
Replace-Text "575÷3=" "256÷6="
Replace-Text "554÷5=" "158÷2="
Replace-Text "405÷4=" "847÷3="
Replace-Text "253÷6=" "844÷6="
Replace-Text "716÷5=" "936÷6="
Replace-Text "370÷4=" "695÷9="
Replace-Text "915÷2=" "670÷2="
Replace-Text "986÷2=" "536÷3="
Replace-Text "844÷5=" "461÷6="
Replace-Text "868÷4=" "414÷6="
Replace-Text "208÷3=" "101÷5="
Replace-Text "763÷2=" "533÷9="
Replace-Text "600÷4=" "750÷3="
Replace-Text "765÷9=" "765÷6="
Replace-Text "197÷4=" "626÷7="
Replace-Text "695÷2=" "519÷3="
Replace-Text "713÷8=" "701÷3="
Replace-Text "886÷8=" "687÷5="
Replace-Text "680÷7=" "464÷6="

Write-Output "Done"
